$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 2907348.78
$ws.Range("C9").Value = 461726.99
$ws.Range("D9").Value = 3369075.77
$ws.Range("E9").Value = 13.70485621342972
$ws.Range("F9").Value = 86.29514378657029
$ws.Range("G9").Value = -55.37638869341093
$ws.Range("H9").Value = -47.49729968585699
$ws.Range("I9").Value = 29020
$ws.Range("J9").Value = 1244
$ws.Range("K9").Value = 30264
$ws.Range("L9").Value = 20860
$ws.Range("M9").Value = 161.5089055608821
$ws.Range("N9").Value = 10.26541939609993
